# PROCESSOR UPDATE 3: UPDATED IMMEDIATE CONTROL SIGNALS
# Rename the "RegDist" control-signal column header to "RegDst", and set the
# RegDst column to 1 for every immediate-type instruction row (the rows where
# ALUSrc == 1), matching the corrected control signal table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller")

# Header rename: "RegDist" -> "RegDst" (column E, row 1)
$ws.Range("E1").Value = "RegDst"

# Correct the RegDst control signal to 1 for immediate instructions
# (addiu, addi, andi, ori, xori, slti, sltiu), which are the rows where
# ALUSrc (column G) is already 1.
$ws.Range("E3").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("E33").Value = 1

# Leave the selection where it ended up after editing the table.
$ws.Range("E33").Select()
